$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The original document wraps the first paragraph ("On Pilgrimage -
# May 1949") in a bookmarkStart/bookmarkEnd pair whose w:id is the
# (schema-invalid, non-numeric) bookmark name itself. Because of that,
# the bookmark never shows up in $d.Bookmarks, so it can't be removed
# through the normal Bookmarks collection API. Deleting paragraph 1
# first collapses both zero-width bookmark anchors down to absolute
# position 0; two subsequent zero-length deletes at position 0 then
# pop them off one at a time, leaving the document bookmark-free -
# matching the target (which drops the bookmark entirely).
# ------------------------------------------------------------------

$d.Paragraphs(1).Range.Delete()
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# Paragraph 1 is now "By Dorothy Day" (bold, no paragraph style).
# Replace it with the two paragraphs pandoc produces for a title
# block: a Title-styled heading (its words/spaces/punctuation each in
# their own run) and an Authors-styled byline (just "Dorothy Day",
# also split word by word, with the "By " prefix and bold formatting
# dropped).
$titleAndAuthor = $d.Paragraphs(1).Range

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">On</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">May</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">1949</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>'

$titleAndAuthor.InsertXML($fragment)
